$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": species counts zeroed out, percentages cleared ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2:B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- Sheet "Species qualification": SoIB 2023 Assessment count zeroed out ---
$wsSpecies = $wb.Worksheets.Item("Species qualification")
$wsSpecies.Range("B5").Value = 0

# --- Sheet "High Priority break-up": recomputed breakdown, now only 2 rows ---
$wsBreakup = $wb.Worksheets.Item("High Priority break-up")
$wsBreakup.Range("A2").Value = "Trend New"
$wsBreakup.Range("B2").Value = 5
$wsBreakup.Range("C2").Value = 21.7
$wsBreakup.Range("D2").Value = 5
$wsBreakup.Range("E2").Value = 21.7

$wsBreakup.Range("A3").Value = "IUCN"
$wsBreakup.Range("B3").Value = 18
$wsBreakup.Range("C3").Value = 78.3
$wsBreakup.Range("D3").Value = 18
$wsBreakup.Range("E3").Value = 78.3

$wsBreakup.Range("A4:E5").ClearContents()
